$d = $word.ActiveDocument
$vt = [char]11

# Title
$d.Paragraphs(1).Range.Text = "LTI Solution"

# Question 1
$d.Paragraphs(2).Range.Text = ("1)what is 2+5?" + $vt + "a)7" + $vt + $vt)

# Question 2
$d.Paragraphs(3).Range.Text = ("2)2" + $vt + "a)3" + $vt + "b)3" + $vt + $vt)

# Question 3
$d.Paragraphs(4).Range.Text = ("3)what ?" + $vt + "a)yes" + $vt + "b)o" + $vt + $vt)
